# Update database and shift the quarterly "cumulative" columns left by one
# period: drop the oldest period (column D) and append a brand-new period
# in column M, mirroring the author's commit ("update database and change
# read_price algorithm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the oldest period column (D). This shifts every later column one
#    slot to the left (values, header labels/shared strings, styles, and
#    column-width bands all move together), exactly like the diff shows.
$ws.Columns.Item(4).Delete()

# 2) Recreate a 13th column (M) with the same per-row formatting as the new
#    last existing column (L), then give it its own width band.
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$ws.Columns.Item(13).ColumnWidth = 28.17

# 3) New period header + publish-date header for column M.
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-27 (2)"

# 4) The publish date that used to read "1401-10-29 (7)" is now corrected
#    to "1402-02-27 (9)" (column I after the shift).
$ws.Range("I9").Value = "1402-02-27 (9)"

# 5) New period's financial figures (column M), row by row.
$ws.Range("M11").Value = 111798328
$ws.Range("M12").Value = -66958743
$ws.Range("M13").Value = 44839585
$ws.Range("M14").Value = -10713596
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 10344702
$ws.Range("M17").Value = 44470691
$ws.Range("M18").Value = 0
$ws.Range("M19").Value = 2137418
$ws.Range("M20").Value = 46608109
$ws.Range("M21").Value = -8134088
$ws.Range("M22").Value = 38474021
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 38474021
$ws.Range("M25").Value = 6412
$ws.Range("M26").Value = 6000000
$ws.Range("M27").Value = 6412

# 6) read_price algorithm change: column I (previously a "-" placeholder in
#    the EPS-after-tax row because that period had zero recorded capital)
#    now resolves to an actual per-share figure.
$ws.Range("I25").Value = 4749
